# Update the "as_of_utc" timestamp column (AA) on the data sheets.
$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-10-30 03:39:30"
$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Range("AA$row").Value = $newTimestamp
    }
}
